# More tests on CS300 - update computed cost values on the COSTS sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Depreciation
$ws.Range("B2").Value = 6999.995406508379
$ws.Range("C2").Value = 1491.7477344706822
$ws.Range("D2").Value = 4.374997129067738
$ws.Range("E2").Value = 2.9166647527118257

# Row 3 - Interest
$ws.Range("B3").Value = 6719.995590248043
$ws.Range("C3").Value = 1432.0778250918547
$ws.Range("D3").Value = 4.199997243905028
$ws.Range("E3").Value = 2.7999981626033517

# Row 4 - Insurance
$ws.Range("B4").Value = 516.2783123628941
$ws.Range("C4").Value = 110.02250117301881
$ws.Range("D4").Value = 0.3226739452268089
$ws.Range("E4").Value = 0.21511596348453926

# Row 5 - DOC Capital
$ws.Range("B5").Value = 14236.269309119318
$ws.Range("C5").Value = 3033.848060735556
$ws.Range("D5").Value = 8.897668318199576
$ws.Range("E5").Value = 5.931778878799717

# Row 7 - Cockpit Crew
$ws.Range("B7").Value = 6757.170232236857
$ws.Range("D7").Value = 4.223231395148037
$ws.Range("E7").Value = 2.815487596765358

# Row 8 - Cabin Crew
$ws.Range("B8").Value = 1266.969418544411
$ws.Range("D8").Value = 0.7918558865902571
$ws.Range("E8").Value = 0.5279039243935046

# Row 9 - DOC Crew
$ws.Range("B9").Value = 8024.139650781269
$ws.Range("D9").Value = 5.015087281738294
$ws.Range("E9").Value = 3.3433915211588627

# Row 11 - Landing charges
$ws.Range("B11").Value = 3452.2578682970907
$ws.Range("C11").Value = 735.7001761819095
$ws.Range("D11").Value = 2.1576611676856823
$ws.Range("E11").Value = 1.4384407784571216

# Row 13 - Ground handling charges
$ws.Range("B13").Value = 508.4273999999998
$ws.Range("C13").Value = 108.34941711356552
$ws.Range("D13").Value = 0.317767125
$ws.Range("E13").Value = 0.21184475

# Row 14 - Noise charges
$ws.Range("B14").Value = 1691.6602304037292
$ws.Range("C14").Value = 360.50456745337505
$ws.Range("D14").Value = 1.057287644002331
$ws.Range("E14").Value = 0.7048584293348874

# Row 15 - Emissions charges
$ws.Range("C15").Value = 282.57982770517066
$ws.Range("D15").Value = 0.82875
$ws.Range("E15").Value = 0.5525

# Row 16 - DOC Charges
$ws.Range("C16").Value = 1.1607354538749344
$ws.Range("D16").Value = 0.0034042044515736314
$ws.Range("E16").Value = 0.002269469634382421

# Row 17 - Airframe Maintenance Charges
$ws.Range("C17").Value = 10.593665665578136
$ws.Range("D17").Value = 0.031069098214285704
$ws.Range("E17").Value = 0.020712732142857137

# Row 18 - Engine Maintenance Charges
$ws.Range("B18").Value = 3395.668913669103
$ws.Range("C18").Value = 723.6406761451127
$ws.Range("D18").Value = 2.12229307104319
$ws.Range("E18").Value = 1.4148620473621265

# Row 20 - DOC Maintenance
$ws.Range("B20").Value = 1609.158191191197
$ws.Range("C20").Value = 342.92280876106537
$ws.Range("D20").Value = 1.0057238694944983
$ws.Range("E20").Value = 0.670482579662999

# Row 21 - Total DOC
$ws.Range("B21").Value = 5937.368638487978
$ws.Range("C21").Value = 1265.2945753288216
$ws.Range("D21").Value = 3.7108553990549877
$ws.Range("E21").Value = 2.473903599369992

# Row 22 - Cash DOC
$ws.Range("B22").Value = 7655.283504533765
$ws.Range("C22").Value = 1631.3941883449377
$ws.Range("D22").Value = 4.784552190333605
$ws.Range("E22").Value = 3.1897014602224028

# Row 25
$ws.Range("B25").Value = 36763.61924640054
$ws.Range("C25").Value = 7834.583101407516
$ws.Range("D25").Value = 22.977262029000347
$ws.Range("E25").Value = 15.318174686000232

# Row 27
$ws.Range("B27").Value = 22527.34993728123
$ws.Range("C27").Value = 4800.735040671961
$ws.Range("D27").Value = 14.079593710800772
$ws.Range("E27").Value = 9.386395807200515
